# Applies the "Updated ptox data and plots" edit:
#  - Appends new PTOX-results rows (41-53) to the "PTOX results" sheet,
#    including one cell with mixed-italic rich text.
#  - Removes the two incomplete trailing rows (28-29) from "Toxindata".
#  - Switches the active/selected sheet & view back to "PTOX results".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. PTOX results — append rows 41-53
# ---------------------------------------------------------------------
$ptox = $wb.Worksheets.Item("PTOX results")

# Register an italic "Calibri 11 / theme color 1" font in the workbook's
# style table (mirrors the font table entry Excel creates the first time
# italics get used) without leaving any extra cell using it directly.
$fontProbe = $ptox.Cells.Item(999, 1)
$fontProbe.Value = "x"
$fontProbe.Font.Italic = $true
$fontProbe.Clear() | Out-Null

# Helper source cells carrying the two formats we need to replicate via
# PasteSpecial (keeps everything on the model's existing style indices
# instead of synthesizing brand-new ones):
#   - date format (numFmtId 14) already used throughout column B
$dateFormatSrc = $ptox.Cells.Item(40, 2)

function Set-Row($r, $a, $b, $c, $d, $e, $f, $g, $h, $cDateStyle) {
    $ptox.Cells.Item($r, 1).Value = $a
    $ptox.Cells.Item($r, 2).Value = $b
    $ptox.Cells.Item($r, 3).Value = $c
    if ($cDateStyle) {
        $dateFormatSrc.Copy() | Out-Null
        $ptox.Cells.Item($r, 3).PasteSpecial(-4122) | Out-Null
    }
    if ($d -ne $null) { $ptox.Cells.Item($r, 4).Value = $d }
    if ($e -ne $null) { $ptox.Cells.Item($r, 5).Value = $e }
    $ptox.Cells.Item($r, 6).Value = $f
    $ptox.Cells.Item($r, 7).Value = $g
    $ptox.Cells.Item($r, 8).Value = $h
    # Dates always pick up column B's existing date display format.
    $dateFormatSrc.Copy() | Out-Null
    $ptox.Cells.Item($r, 2).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

Set-Row 41 "FAL" 44859 "FAL_SG_221025" "NA" $null "N" "Greenwater" "PTOX cyanobacteria were not observed." $false
Set-Row 42 "FAL" 44859 "DUP_SG_221025" "NA" $null "N" "Greenwater" "PTOX cyanobacteria were not observed." $false
Set-Row 43 "HOL" 44847 "HOL_SG_221013" "Phormidium" 1 "N" "Greenwater" "The PTOX cyanobacterium cf. Phormidium sp. (~1 filament per mL) was observed." $false
Set-Row 44 "HLT" 44847 "HLT_SG_221013" "NA" $null "N" "Greenwater" "PTOX cyanobacteria were not observed." $false
Set-Row 45 "HLT" 44847 "DUP_SG_221013" "NA" $null "N" "Greenwater" "PTOX cyanobacteria were not observed." $false
Set-Row 46 "FAL" 44830 "FAL_SG_220926" "NA" $null "N" "Greenwater" "PTOX cyanobacteria were not observed." $true

# Row 47: first half of the mixed-species note (Microcystis + Aphanizomenon
# observed together) - build the rich-text note once, with the two genus
# names italicised, then reuse the exact same note for row 48.
Set-Row 47 "FAL" 44830 "DUP_SG_220926" "Microcystis" 1 "N" "Greenwater" "The PTOX cyanobacteria Microcystis sp. (1 colonly per mL) and Aphanizomenon sp. (2 filaments per mL) were observed." $true
$noteCell = $ptox.Cells.Item(47, 8)
$noteCell.Characters(24, 11).Font.Italic = $true   # "Microcystis"
$noteCell.Characters(35, 28).Font.Italic = $false  # " sp. (1 colonly per mL) and "
$noteCell.Characters(63, 14).Font.Italic = $true   # "Aphanizomenon "
$noteCell.Characters(77, 39).Font.Italic = $false  # "sp. (2 filaments per mL) were observed."

Set-Row 48 "FAL" 44830 "DUP_SG_220926" "Aphanizomenon" 2 "N" "Greenwater" $null $true
$noteCell.Copy() | Out-Null
$ptox.Cells.Item(48, 8).PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

Set-Row 49 "HOL" 44812 "HOL_SG_220908" "NA" $null "N" "Greenwater" "PTOX cyanobacteria were not observed." $true
Set-Row 50 "HLT" 44812 "HLT_SG_220908" "NA" $null "N" "Greenwater" "PTOX cyanobacteria were not observed." $true
Set-Row 51 "HLT" 44812 "DUP_SG_220908" "Planktothrix" 3 "N" "Greenwater" "The PTOX cyanobacteria Planktothrix sp. (3 filaments per mL) and 1 nostocalean filament per mL were observed. Specialized cells (i.e. akinetes, heterocytes) were not present on the nostocalean filament, preventing genus level identification." $true
Set-Row 52 "FAL" 44804 "FAL_SG_2020831" "Planktothrix" 1 "N" "Greenwater" "The PTOX cyanobacterium Planktothrix sp. (1 filament per mL) was observed." $true
Set-Row 53 "FAL" 44804 "DUP_SG_2020831" "Planktothrix" 1 "N" "Greenwater" "The PTOX cyanobacterium Planktothrix sp. (1 filament per mL) was observed." $true

# Print area / page setup picked up when the sheet was next printed/previewed.
$ptox.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 2. Toxindata — drop the two incomplete trailing rows (28 & 29)
# ---------------------------------------------------------------------
$tox = $wb.Worksheets.Item("Toxindata")
$tox.Rows.Item(28).Resize(2).Delete() | Out-Null

# ---------------------------------------------------------------------
# 3. View state — PTOX results becomes the active/selected sheet again
# ---------------------------------------------------------------------
$ptox.Activate()
$ptox.Range("H27").Select() | Out-Null

Write-Output "done"
